$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.546.19"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.468.33"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9599"
$ws.Range("E5").Value = "  +5.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.67"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3554"
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.087"
$ws.Range("E9").Value = "  +6.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.44"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06615"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.456"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("E14").Value = "  +3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.156"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9605"
$ws.Range("E16").Value = "  +3.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001019"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.470.13"
$ws.Range("E18").Value = "  +2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05957"
$ws.Range("E19").Value = "  +5.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.67"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.472"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("E22").Value = "  +2.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.25"
$ws.Range("E23").Value = "  +4.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.261"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.534.22"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.73"
$ws.Range("E26").Value = "  +5.28%  "
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.07"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.631.60"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "113.73"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.850"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07941"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.905"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7962"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.245"
$ws.Range("E35").Value = "  +10.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.458"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05751"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.700"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9611"
$ws.Range("E39").Value = "  +3.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02028"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.31"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1863"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.281"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5244"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.513"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.05"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "119.10"
$ws.Range("E47").Value = "  +2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5172"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.799"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06436"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9932"
$ws.Range("E51").Value = "  +1.52%  "
